# Update "想去人数" (F column) counts for several events on both the
# "展览" and "全部类型" sheets, per the source diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 533
$ws1.Range("F5").Value = 507
$ws1.Range("F8").Value = 22
$ws1.Range("F14").Value = 822
$ws1.Range("F15").Value = 279
$ws1.Range("F16").Value = 581
$ws1.Range("F18").Value = 1326
$ws1.Range("F21").Value = 1181
$ws1.Range("F22").Value = 2860
$ws1.Range("F24").Value = 700
$ws1.Range("F30").Value = 3066
$ws1.Range("F31").Value = 594
$ws1.Range("F33").Value = 1390

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 533
$ws4.Range("F7").Value = 507
$ws4.Range("F15").Value = 22
$ws4.Range("F26").Value = 822
$ws4.Range("F27").Value = 279
$ws4.Range("F28").Value = 581
$ws4.Range("F30").Value = 1326
$ws4.Range("F33").Value = 1181
$ws4.Range("F34").Value = 2860
$ws4.Range("F36").Value = 700
$ws4.Range("F44").Value = 3066
$ws4.Range("F45").Value = 594
$ws4.Range("F47").Value = 1390
